$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the dSF (column F) values for the re-pulled rows.
$ws.Range("F2").Value = 0
$ws.Range("F3").Value = -1
$ws.Range("F4").Value = 2
$ws.Range("F8").Value = -6
$ws.Range("F11").Value = 0
$ws.Range("F12").Value = -1
